# player can collect items and drop to the map
#
# Applies the authored changes:
#  - "TODO's" sheet: three new backlog rows (Columns randomness / wall color /
#    wall traps) with matching status + detail notes, taller rows for the
#    wrapped note cells, and the selection left on C14.
#  - "Logs" sheet: a new dated entry describing the item drop/collect work,
#    becomes the active/selected sheet, selection left on B57.
#  - "TODO Before 0.0.1" sheet: scrolled down a bit (view-only).
#  - "insane calcs" sheet: no longer the active tab (Logs takes over).

$wb = $excel.ActiveWorkbook

$wsTodo   = $wb.Worksheets.Item("TODO's")
$wsLogs   = $wb.Worksheets.Item("Logs")
$wsBefore = $wb.Worksheets.Item("TODO Before 0.0.1")

# ---------------------------------------------------------------------
# "TODO Before 0.0.1" sheet - scroll the view down a few rows (the cell
# selection itself, A47, is untouched by the authored edit - only the
# window's scroll/top-left-row changes).
# (done first so the later "Logs" activation below wins and ends up as
# the workbook's active tab, matching the authored edit)
# ---------------------------------------------------------------------
$wsBefore.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# "TODO's" sheet - add rows 12-14 (Nr 11-13) under the existing backlog.
# Values are written in the same order the author must have typed them in
# (note first, then the three todo descriptions, then the second note) so
# new shared strings land on the same indices as the source file.
# ---------------------------------------------------------------------
$wsTodo.Range("E12").Value = "should be different and sometimes not collumns but some strange objects"
$wsTodo.Range("B12").Value = "Columns randomnless"
$wsTodo.Range("B13").Value = "wall color - new levels or rooms should looks different little bit"
$wsTodo.Range("B14").Value = "in wall could be traps "
$wsTodo.Range("E14").Value = "when player commes to close to wall - wall openes and some enemy shoot player"

$wsTodo.Range("A12").Value = 11
$wsTodo.Range("A13").Value = 12
$wsTodo.Range("A14").Value = 13

# status column - copy formatting from the row above, then set the same
# "todo" text already used by the rest of the sheet.
$wsTodo.Range("C11").Copy()
$wsTodo.Range("C12").PasteSpecial(-4122)
$wsTodo.Range("C12").Value = "todo"

$wsTodo.Range("C11").Copy()
$wsTodo.Range("C13").PasteSpecial(-4122)
$wsTodo.Range("C13").Value = "todo"

$wsTodo.Range("C11").Copy()
$wsTodo.Range("C14").PasteSpecial(-4122)
$wsTodo.Range("C14").Value = "todo"

# rows with a note in column E wrap to two lines, same as row 6 above.
$wsTodo.Rows.Item(12).RowHeight = 30
$wsTodo.Rows.Item(14).RowHeight = 30

$wsTodo.Range("C14").Select()

# ---------------------------------------------------------------------
# "Logs" sheet - append the new dev-log entry and make it the active tab.
# ---------------------------------------------------------------------
$wsLogs.Range("A55").Copy()
$wsLogs.Range("A56").PasteSpecial(-4122)
$wsLogs.Range("A56").Value = 45566

$wsLogs.Range("B55").Copy()
$wsLogs.Range("B56").PasteSpecial(-4122)
$wsLogs.Range("B56").Value = "items can be droped and collected to invenory!"

$wsLogs.Activate()
$wsLogs.Range("B57").Select()
